# Update the "data updated" timestamp in A1 (10:52 -> 11:22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 11:22"

# --- Row 13 (Rusia) - refreshed case numbers ---
$ws.Range("B13").Value = 80949
$ws.Range("C13").Value = 6361
$ws.Range("D13").Value = 6767
$ws.Range("E13").Value = 73435
$ws.Range("F13").Value = 2300
$ws.Range("G13").Value = 66
$ws.Range("H13").Value = 747

# --- Rows 15/16: Belgica moves above Canada (both refreshed + swapped) ---
# Row 15 becomes Belgica with newly refreshed figures
$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 46134
$ws.Range("C15").Value = 809
$ws.Range("D15").Value = 10785
$ws.Range("E15").Value = 28255
$ws.Range("F15").Value = 891
$ws.Range("G15").Value = 177
$ws.Range("H15").Value = 7094

# Row 16 becomes Canada, carrying the figures Canada previously had in row 15
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 45354
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 16425
$ws.Range("E16").Value = 26464
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 2465

# --- Row 27 (Austria) - refreshed case numbers ---
$ws.Range("B27").Value = 15225
$ws.Range("C27").Value = 77
$ws.Range("D27").Value = 12282
$ws.Range("E27").Value = 2401
$ws.Range("F27").Value = 145
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 542

# --- Rows 39/40: Indonesia moves above Ucrania (both refreshed + swapped) ---
# Row 39 becomes Indonesia with newly refreshed figures
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 8882
$ws.Range("C39").Value = 275
$ws.Range("D39").Value = 1107
$ws.Range("E39").Value = 7032
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 23
$ws.Range("H39").Value = 743

# Row 40 becomes Ucrania, carrying the figures Ucrania previously had in row 39
$ws.Range("A40").Value = "Ucrania"
$ws.Range("B40").Value = 8617
$ws.Range("C40").Value = 492
$ws.Range("D40").Value = 840
$ws.Range("E40").Value = 7568
$ws.Range("F40").Value = 107
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 209

# --- Row 138 (Brunei) - refreshed active/recovered counts ---
$ws.Range("D138").Value = 123
$ws.Range("E138").Value = 14

# --- Row 188 (Namibia) - refreshed active/recovered counts ---
$ws.Range("D188").Value = 8
$ws.Range("E188").Value = 8
